$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# DOMA-746: rename the `d.ticket[...]` export-mapper placeholders to
# `d.tickets[...]` (row 2 = item i, row 3 = item i + 1).
$ws.Range("A2").Value = "{d.tickets[i].address}"
$ws.Range("B2").Value = "{d.tickets[i].processing}"
$ws.Range("C2").Value = "{d.tickets[i].completed}"
$ws.Range("D2").Value = "{d.tickets[i].canceled}"
$ws.Range("E2").Value = "{d.tickets[i].deferred}"
$ws.Range("F2").Value = "{d.tickets[i].closed}"
$ws.Range("G2").Value = "{d.tickets[i].new_or_reopened}"

$ws.Range("A3").Value = "{d.tickets[i + 1].address}"
$ws.Range("B3").Value = "{d.tickets[i + 1].processing}"
$ws.Range("C3").Value = "{d.tickets[i + 1].completed}"
$ws.Range("D3").Value = "{d.tickets[i + 1].canceled}"
$ws.Range("E3").Value = "{d.tickets[i + 1].deferred}"
$ws.Range("F3").Value = "{d.tickets[i + 1].closed}"
$ws.Range("G3").Value = "{d.tickets[i + 1].new_or_reopened}"

# Move the saved cursor/selection to D28, matching the workbook's stored view state.
$ws.Range("D28").Select()
